$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Problem 1 results data (rows 6-15, columns E-H)
$ws.Range("G6").Value = 56
$ws.Range("H6").Value = 180

$ws.Range("E7").Value = 2.605
$ws.Range("F7").Value = 1458
$ws.Range("G7").Value = 1459
$ws.Range("H7").Value = 5960

$ws.Range("E8").Value = 0.037
$ws.Range("F8").Value = 21
$ws.Range("G8").Value = 22
$ws.Range("H8").Value = 84

$ws.Range("E9").Value = 0.19900000000000001
$ws.Range("F9").Value = 101
$ws.Range("G9").Value = 271
$ws.Range("H9").Value = 414

$ws.Range("E10").Value = 0.092
$ws.Range("F10").Value = 55
$ws.Range("G10").Value = 57
$ws.Range("H10").Value = 224

$ws.Range("E11").Value = 5.23
$ws.Range("F11").Value = 4229
$ws.Range("G11").Value = 4230
$ws.Range("H11").Value = 17023

$ws.Range("E12").Value = 0.006
$ws.Range("F12").Value = 7
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = 28

$ws.Range("E13").Value = 0.052
$ws.Range("F13").Value = 55
$ws.Range("G13").Value = 57
$ws.Range("H13").Value = 224

$ws.Range("E14").Value = 0.061
$ws.Range("F14").Value = 41
$ws.Range("G14").Value = 43
$ws.Range("H14").Value = 170

$ws.Range("E15").Value = 0.874
$ws.Range("F15").Value = 11
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = 50

# Update the selection / view to match the committed state
$ws.Range("D3:I15").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
